# Data-cleaning pass: correct the truncated "SELECT" board-result codes
# (and all strings that contain "SELECT" as a substring, e.g. "NON-SELECT",
# "NONSELECT", "BOARD NON-SELECT") to their fully spelled "SELECTED" form
# across the worksheet, matching the commit "add additional cleaning function".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPart = 2 (allow substring matches so "NON-SELECT" -> "NON-SELECTED"),
# xlByRows = 1, MatchCase = False, MatchByte = False
$ws.Cells.Replace("SELECT", "SELECTED", 2, 1, $false, $false, $false, $false)
